$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF column (58) holds the "Date" field; every data row (2-31) has the
# wrong value "5-12-2007-08" and needs to become "2008-05-12".
# Force the range to Text format first so Excel doesn't auto-convert the
# date-shaped string into a date serial number.
$range = $ws.Range("BF2:BF31")
$range.NumberFormat = "@"

for ($r = 2; $r -le 31; $r++) {
    $ws.Cells.Item($r, 58).Value = "2008-05-12"
}
